# Updated symbol list on Thu Dec 15 07:38:40 UTC 2022 with GitHub Actions
#
# Applies the latest price/volume-label refresh to the crypto tracker sheet.
# Numeric-looking values are written as literal text (matching the sheet's
# existing inline-string cells), so each cell is briefly switched to Text
# format before the assignment and then restored to the Normal style so no
# stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Value
    )
    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Value
    $cell.Style = "Normal"
}

Set-TextValue "D2"  "265.02"
Set-TextValue "D3"  "22.73"
Set-TextValue "D4"  "6.270"
Set-TextValue "D5"  "0.06160"
Set-TextValue "D6"  "3.592"
Set-TextValue "D7"  "6.714"
Set-TextValue "D8"  "1.355"
Set-TextValue "D9"  "0.8283"
Set-TextValue "D11" "0.1615"
Set-TextValue "D12" "0.08200"
Set-TextValue "D13" "0.03393"
Set-TextValue "D14" "0.03139"
Set-TextValue "D15" "0.09249"
Set-TextValue "D16" "3.920"
Set-TextValue "D17" "0.001704"
Set-TextValue "D18" "0.04796"
Set-TextValue "D19" "0.006280"
Set-TextValue "D20" "0.005926"
Set-TextValue "D21" "0.001101"
Set-TextValue "D23" "3.768"
Set-TextValue "D24" "2.301"
Set-TextValue "D26" "0.1245"
Set-TextValue "D40" "0.04649"
Set-TextValue "D41" "0.006951"
Set-TextValue "D43" "0.003400"
Set-TextValue "D44" "0.01039"
$ws.Range("E44").Value = "43LocalTradersLCT"
Set-TextValue "D45" "0.00006157"
Set-TextValue "D47" "0.7781"
Set-TextValue "D48" "0.2036"
Set-TextValue "D49" "0.00001400"
$ws.Range("E49").Value = "48CryptobidCoinCBCWorstin24h"
